$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 257; existing rows 257-322 shift down to 258-323.
$ws.Rows.Item(257).Insert()

# Populate the newly inserted row 257 with the new weekly data point.
$ws.Cells.Item(257, 1).Value2  = 11
$ws.Cells.Item(257, 2).Value2  = "Vega Monumental Concepción"
$ws.Cells.Item(257, 3).Value2  = "Bíobío"
$ws.Cells.Item(257, 4).Value2  = 44642
$ws.Cells.Item(257, 5).Value2  = 8
$ws.Cells.Item(257, 6).Value2  = "Fruta"
$ws.Cells.Item(257, 7).Value2  = 100101
$ws.Cells.Item(257, 8).Value2  = "Berries"
$ws.Cells.Item(257, 9).Value2  = 100112025
$ws.Cells.Item(257, 10).Value2 = "Frutilla"
$ws.Cells.Item(257, 11).Value2 = "Sin especificar"
$ws.Cells.Item(257, 12).Value2 = "Primera"
$ws.Cells.Item(257, 13).Value2 = 200
$ws.Cells.Item(257, 14).Value2 = 7000
$ws.Cells.Item(257, 15).Value2 = 7500
$ws.Cells.Item(257, 16).Value2 = 7300
$ws.Cells.Item(257, 17).Value2 = "$/bandeja 7 kilos"
$ws.Cells.Item(257, 18).Value2 = "Región del Maule"
$ws.Cells.Item(257, 19).Value2 = 1043
$ws.Cells.Item(257, 20).Value2 = 7
